$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws "D2" "62.341.12"
Set-TextCell $ws "E2" "  +0.04%  "
Set-TextCell $ws "D3" "2.465.26"
Set-TextCell $ws "E3" "  +1.68%  "
Set-TextCell $ws "E4" "  +0.05%  "
Set-TextCell $ws "D5" "582.58"
Set-TextCell $ws "E5" "  +0.75%  "
Set-TextCell $ws "E6" "  -0.07%  "
Set-TextCell $ws "E7" "  +0.04%  "
Set-TextCell $ws "E8" "  +1.36%  "
Set-TextCell $ws "D9" "2.462.44"
Set-TextCell $ws "E9" "  +1.73%  "
Set-TextCell $ws "E10" "  +4.09%  "
Set-TextCell $ws "E11" "  +2.73%  "
Set-TextCell $ws "E12" "  +0.22%  "
Set-TextCell $ws "E13" "  -0.59%  "
Set-TextCell $ws "E14" "  -0.47%  "
Set-TextCell $ws "E15" "  +0.62%  "
Set-TextCell $ws "D16" "2.901.75"
Set-TextCell $ws "D17" "62.208.43"
Set-TextCell $ws "E17" "  +0.20%  "
Set-TextCell $ws "D18" "2.464.47"
Set-TextCell $ws "E18" "  +1.99%  "
Set-TextCell $ws "E19" "  -2.27%  "
Set-TextCell $ws "D20" "7.35"
Set-TextCell $ws "E20" "  +3.91%  "
Set-TextCell $ws "D21" "327.45"
Set-TextCell $ws "E21" "  -0.65%  "
Set-TextCell $ws "E22" "  -0.33%  "
Set-TextCell $ws "D23" "6.36"
Set-TextCell $ws "E23" "  +6.55%  "
Set-TextCell $ws "E24" "  -0.48%  "
Set-TextCell $ws "E25" "  -0.02%  "
Set-TextCell $ws "D26" "65.48"
Set-TextCell $ws "E26" "  +0.00%  "
Set-TextCell $ws "D27" "9.12"
Set-TextCell $ws "E27" "  +0.48%  "
Set-TextCell $ws "D28" "590.28"
Set-TextCell $ws "E28" "  -6.53%  "
Set-TextCell $ws "D29" "2.576.58"
Set-TextCell $ws "E29" "  +1.41%  "
Set-TextCell $ws "B30" "PEPE"
Set-TextCell $ws "C30" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell $ws "D30" "0.0₃0950"
Set-TextCell $ws "E30" "  +0.33%  "
Set-TextCell $ws "B31" "Binance-PegBSC-USD"
Set-TextCell $ws "C31" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextCell $ws "D31" "0.999"
Set-TextCell $ws "E31" "  -0.27%  "
Set-TextCell $ws "D32" "8.00"
Set-TextCell $ws "E32" "  -0.25%  "
Set-TextCell $ws "E33" "  -2.70%  "
Set-TextCell $ws "E34" "  +0.33%  "
Set-TextCell $ws "E35" "  -3.12%  "
Set-TextCell $ws "E36" "  -0.11%  "
Set-TextCell $ws "D37" "4.80"
Set-TextCell $ws "E37" "  -2.93%  "
Set-TextCell $ws "E38" "  -2.10%  "
Set-TextCell $ws "E39" "  +0.01%  "
Set-TextCell $ws "D40" "151.50"
Set-TextCell $ws "E40" "  +2.23%  "
Set-TextCell $ws "E41" "  +0.06%  "
Set-TextCell $ws "E42" "  -0.64%  "
Set-TextCell $ws "E43" "  -1.37%  "
Set-TextCell $ws "D44" "42.35"
Set-TextCell $ws "E44" "  +0.33%  "
Set-TextCell $ws "E45" "  +0.00%  "
Set-TextCell $ws "D46" "2.43"
Set-TextCell $ws "E46" "  -1.90%  "
Set-TextCell $ws "D47" "0.0₆0289"
Set-TextCell $ws "E47" "  +21.83%  "
Set-TextCell $ws "D48" "143.37"
Set-TextCell $ws "E48" "  +0.08%  "
Set-TextCell $ws "E49" "  -1.45%  "
Set-TextCell $ws "E50" "  +1.79%  "
Set-TextCell $ws "D51" "19.99"
Set-TextCell $ws "E51" "  +2.83%  "
